# Add "Race Code (Adams County)" / PersonRaceCode row to the Booking Report
# sheet, inserted as the new row 22 (pushing the former rows 22+ down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 22 (the "SSN" row), copying the
# formatting of the row above (row 21) as Excel normally does on insert.
$ws.Rows.Item(22).Insert()

# Populate the new row with the Race Code (Adams County) mapping.
$ws.Range("A22").Value = "x-ext"
$ws.Range("B22").Value = "Race Code (Adams County)"
$ws.Range("C22").Value = "Person Race"
$ws.Range("E22").Value = "/br-doc:BookingReport/nc:Person[@structures:id=/br-doc:BookingReport/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/ac-bkg-codes:PersonRaceCode"

# Match the row height used for this kind of wrapped, multi-line entry.
$ws.Rows.Item(22).RowHeight = 45

# Restore the selection to the cell the author left active.
$ws.Activate()
$ws.Range("D20").Select()
